# Atualização automática de preços de eletricidade
# Updates row 2 of the SpotPTTable with the latest daily hourly spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45964

$ws.Range("B2").Value = 88.95999999999999
$ws.Range("C2").Value = 83.81999999999999
$ws.Range("D2").Value = 83
$ws.Range("E2").Value = 82.56
$ws.Range("F2").Value = 82.81999999999999
$ws.Range("G2").Value = 86.01000000000001
$ws.Range("H2").Value = 103.21
$ws.Range("I2").Value = 133.09
$ws.Range("J2").Value = 113.73
$ws.Range("K2").Value = 70.89
$ws.Range("L2").Value = 34.14
$ws.Range("M2").Value = 19.97
$ws.Range("N2").Value = 16.49
$ws.Range("O2").Value = 13.72
$ws.Range("P2").Value = 18.62
$ws.Range("Q2").Value = 36.29
$ws.Range("R2").Value = 63
$ws.Range("S2").Value = 110.37
$ws.Range("T2").Value = 125.03
$ws.Range("U2").Value = 147.6
$ws.Range("V2").Value = 151.42
$ws.Range("W2").Value = 130.49
$ws.Range("X2").Value = 102.28
$ws.Range("Y2").Value = 91.84999999999999
$ws.Range("Z2").Value = 82.89

$ws.Range("AB2").Value = 119.01
$ws.Range("AD2").Value = 140.95
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 136.32
$ws.Range("AG2").Value = "3h-16h"
